# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Thu Jul 18 15:29:18 UTC 2024 with GitHub Actions".
# Column D (Price) and column E (Volume(1h)) values are updated in place for
# rows 2-51 of the active sheet.
#
# Price values that look like a plain decimal number (e.g. "7.15") would be
# auto-coerced from text to a number by a plain .Value assignment (exactly like
# typing them into Excel would), which loses formatting such as trailing zeros
# ("0.0730" -> 0.073). To keep them as literal text we write them with a
# leading quote-prefix (forcing text) and then restore the "Normal" cell style
# so no stray number-format/style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.590.79"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "3.404.49"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Formula = "'567.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").Formula = "'156.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.402.71"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Formula = "'0.571"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.49%  "
$ws.Range("D10").Formula = "'7.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Formula = "'0.118"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.30%  "
$ws.Range("E12").Value = "  -5.61%  "
$ws.Range("D13").Value = "3.992.65"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Formula = "'26.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.78%  "
$ws.Range("D16").Formula = "'0.0000173"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.82%  "
$ws.Range("D17").Value = "63.672.15"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "3.412.99"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Formula = "'6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.67%  "
$ws.Range("D20").Formula = "'13.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.94%  "
$ws.Range("D21").Formula = "'374.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Formula = "'7.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Formula = "'71.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").Formula = "'0.519"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.98%  "
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Formula = "'9.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.90%  "
$ws.Range("D28").Formula = "'0.176"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Formula = "'5.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("D31").Formula = "'1.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.13%  "
$ws.Range("D32").Formula = "'1.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("D33").Formula = "'22.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Formula = "'6.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.22%  "
$ws.Range("D35").Formula = "'1.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.18%  "
$ws.Range("D36").Formula = "'160.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Formula = "'1.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("D38").Formula = "'0.815"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").Formula = "'26.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("D40").Formula = "'0.0730"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.12%  "
$ws.Range("D41").Value = "2.781.14"
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("D42").Formula = "'42.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").Formula = "'4.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.07%  "
$ws.Range("D44").Formula = "'6.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.00%  "
$ws.Range("E45").Value = "  -5.41%  "
$ws.Range("D46").Formula = "'25.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("E47").Value = "  +7.11%  "
$ws.Range("D48").Formula = "'325.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  -6.24%  "
$ws.Range("E50").Value = "  -4.43%  "
$ws.Range("D51").Formula = "'0.822"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.64%  "
